$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111756210
$ws.Range("B2").Value = 56398
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 454473
$ws.Range("R2").Value = 7077762
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").Value = "ringhack"

# Row 3
$ws.Range("A3").Value = 111756218
$ws.Range("B3").Value = 56398
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("Q3").Value = 454467
$ws.Range("R3").Value = 7077694
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").Value = "ringhack äldre"

# Row 4
$ws.Range("A4").Value = 111756219
$ws.Range("B4").Value = 56398
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 454504
$ws.Range("R4").Value = 7077712
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").Value = "ringhack äldre"

# Row 5
$ws.Range("A5").Value = 111756215
$ws.Range("B5").Value = 56398
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("Q5").Value = 454579
$ws.Range("R5").Value = 7077736
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("AC5").Value = "ringhack äldre"

# Row 6
$ws.Range("A6").Value = 111756216
$ws.Range("B6").Value = 56398
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("Q6").Value = 454531
$ws.Range("R6").Value = 7077805
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").Value = "ringhack äldre"

# Row 7
$ws.Range("A7").Value = 111756232
$ws.Range("B7").Value = 89423
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 5432
$ws.Range("F7").Value = "Granticka"
$ws.Range("G7").Value = "Porodaedalea chrysoloma"
$ws.Range("H7").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q7").Value = 454523
$ws.Range("R7").Value = 7077787
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AC7").ClearContents()

# Row 8
$ws.Range("A8").Value = 111756175
$ws.Range("B8").Value = 89405
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1202
$ws.Range("F8").Value = "Ullticka"
$ws.Range("G8").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H8").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q8").Value = 454431
$ws.Range("R8").Value = 7077752
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").ClearContents()

# Row 9
$ws.Range("A9").Value = 111756213
$ws.Range("B9").Value = 56398
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("Q9").Value = 454622
$ws.Range("R9").Value = 7077658
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()
$ws.Range("AC9").Value = "ringhack äldre"

# Row 10
$ws.Range("A10").Value = 111756231
$ws.Range("B10").Value = 89423
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 5432
$ws.Range("F10").Value = "Granticka"
$ws.Range("G10").Value = "Porodaedalea chrysoloma"
$ws.Range("H10").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q10").Value = 454535
$ws.Range("R10").Value = 7077734
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()
$ws.Range("AC10").ClearContents()

# Row 11
$ws.Range("A11").Value = 111756179
$ws.Range("B11").Value = 90087
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 3298
$ws.Range("F11").Value = "Trådticka"
$ws.Range("G11").Value = "Climacocystis borealis"
$ws.Range("H11").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q11").Value = 454498
$ws.Range("R11").Value = 7077773
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()
$ws.Range("AC11").ClearContents()

# Row 12
$ws.Range("A12").Value = 111756211
$ws.Range("B12").Value = 56398
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("Q12").Value = 454467
$ws.Range("R12").Value = 7077794
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()
$ws.Range("AC12").Value = "ringhack"

# Row 13
$ws.Range("A13").Value = 111756233
$ws.Range("B13").Value = 89423
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 5432
$ws.Range("F13").Value = "Granticka"
$ws.Range("G13").Value = "Porodaedalea chrysoloma"
$ws.Range("H13").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q13").Value = 454467
$ws.Range("R13").Value = 7077820
$ws.Range("Z13").ClearContents()
$ws.Range("AB13").ClearContents()
$ws.Range("AC13").ClearContents()

# Row 14
$ws.Range("A14").Value = 111756217
$ws.Range("B14").Value = 56398
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 100109
$ws.Range("F14").Value = "Tretåig hackspett"
$ws.Range("G14").Value = "Picoides tridactylus"
$ws.Range("H14").Value = "(Linnaeus, 1758)"
$ws.Range("Q14").Value = 454478
$ws.Range("R14").Value = 7077819
$ws.Range("Z14").ClearContents()
$ws.Range("AB14").ClearContents()
$ws.Range("AC14").Value = "ringhack"

# Row 15
$ws.Range("A15").Value = 111756209
$ws.Range("B15").Value = 56398
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 100109
$ws.Range("F15").Value = "Tretåig hackspett"
$ws.Range("G15").Value = "Picoides tridactylus"
$ws.Range("H15").Value = "(Linnaeus, 1758)"
$ws.Range("Q15").Value = 454565
$ws.Range("R15").Value = 7077682
$ws.Range("Z15").ClearContents()
$ws.Range("AB15").ClearContents()
$ws.Range("AC15").Value = "ringhack"

# Row 16
$ws.Range("A16").Value = 111756176
$ws.Range("B16").Value = 89405
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 1202
$ws.Range("F16").Value = "Ullticka"
$ws.Range("G16").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H16").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q16").Value = 454476
$ws.Range("R16").Value = 7077724
$ws.Range("Z16").ClearContents()
$ws.Range("AB16").ClearContents()
$ws.Range("AC16").ClearContents()

# Row 17
$ws.Range("A17").Value = 111756181
$ws.Range("B17").Value = 90087
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 3298
$ws.Range("F17").Value = "Trådticka"
$ws.Range("G17").Value = "Climacocystis borealis"
$ws.Range("H17").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q17").Value = 454460
$ws.Range("R17").Value = 7077741
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()
$ws.Range("AC17").ClearContents()

# Row 18
$ws.Range("A18").Value = 111756214
$ws.Range("B18").Value = 56398
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("Q18").Value = 454612
$ws.Range("R18").Value = 7077689
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
$ws.Range("AC18").Value = "ringhack"

# Row 19
$ws.Range("A19").Value = 111756180
$ws.Range("B19").Value = 90087
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 3298
$ws.Range("F19").Value = "Trådticka"
$ws.Range("G19").Value = "Climacocystis borealis"
$ws.Range("H19").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q19").Value = 454561
$ws.Range("R19").Value = 7077725
$ws.Range("Z19").ClearContents()
$ws.Range("AB19").ClearContents()
$ws.Range("AC19").ClearContents()

# Row 20
$ws.Range("A20").Value = 111756220
$ws.Range("B20").Value = 56398
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("Q20").Value = 454526
$ws.Range("R20").Value = 7077711
$ws.Range("Z20").ClearContents()
$ws.Range("AB20").ClearContents()
$ws.Range("AC20").Value = "ringhack"
